$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '42.404.52'
$ws.Range('D3').Value = '2.219.44'
$ws.Range('E3').Value = '  -2.20%  '
$ws.Range('E4').Value = '  +0.31%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '109.45'
$ws.Range('E5').Value = '  -8.37%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '289.30'
$ws.Range('E6').Value = '  +6.81%  '
$ws.Range('E7').Value = '  -3.35%  '
$ws.Range('E8').Value = '  -0.27%  '
$ws.Range('E9').Value = '  -4.55%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '42.95'
$ws.Range('E10').Value = '  -9.54%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '54.18'
$ws.Range('E12').Value = '  -0.14%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '8.52'
$ws.Range('E13').Value = '  -8.88%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.991'
$ws.Range('E14').Value = '  +8.17%  '
$ws.Range('E15').Value = '  -3.17%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '14.77'
$ws.Range('E16').Value = '  -6.38%  '
$ws.Range('D17').Value = '2.550.10'
$ws.Range('E17').Value = '  -2.34%  '
$ws.Range('D18').Value = '2.235.74'
$ws.Range('E18').Value = '  -1.33%  '
$ws.Range('D19').Value = '42.252.13'
$ws.Range('E19').Value = '  -3.08%  '
$ws.Range('E20').Value = '  -5.22%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '7.01'
$ws.Range('E21').Value = '  +1.38%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '72.53'
$ws.Range('E22').Value = '  +0.08%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '3.36'
$ws.Range('E23').Value = '  +11.80%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.36'
$ws.Range('E24').Value = '  -1.41%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '228.99'
$ws.Range('E25').Value = '  -2.45%  '
$ws.Range('E26').Value = '  -8.10%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.999'
$ws.Range('E27').Value = '  -1.88%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '11.29'
$ws.Range('E28').Value = '  -7.83%  '
$ws.Range('E29').Value = '  -2.57%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '172.96'
$ws.Range('E30').Value = '  -1.03%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '36.93'
$ws.Range('E31').Value = '  -11.26%  '
$ws.Range('E32').Value = '  -7.04%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '20.67'
$ws.Range('E33').Value = '  -3.97%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.0871'
$ws.Range('E34').Value = '  -5.20%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '5.54'
$ws.Range('E35').Value = '  -3.16%  '
$ws.Range('E36').Value = '  +4.52%  '
$ws.Range('E37').Value = '  -4.45%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '4.08'
$ws.Range('E38').Value = '  -5.50%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.0364'
$ws.Range('E39').Value = '  -3.90%  '
$ws.Range('E40').Value = '  -5.79%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '72.80'
$ws.Range('E41').Value = '  +0.92%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '2.36'
$ws.Range('E42').Value = '  -7.77%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.227'
$ws.Range('E43').Value = '  -5.04%  '
$ws.Range('E44').Value = '  -0.16%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '12.14'
$ws.Range('E45').Value = '  -11.15%  '
$ws.Range('E46').Value = '  -6.27%  '
$ws.Range('E47').Value = '  -7.58%  '
$ws.Range('E48').Value = '  +5.57%  '
$ws.Range('E49').Value = '  -2.15%  '
$ws.Range('B50').Value = 'FraxShare'
$ws.Range('C50').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '8.35'
$ws.Range('E50').Value = '  -2.92%  '
$ws.Range('B51').Value = 'Aave'
$ws.Range('C51').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '100.32'
$ws.Range('E51').Value = '  -2.50%  '
